$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-12-09", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-10", 2)
